# 2017-07-30-Feeding.xlsx update: add 8/2 data
$wb = $excel.ActiveWorkbook

# --- Rename the first sheet, matching the workbook.xml diff ---
$wsMorning = $wb.Worksheets.Item("morning-concentrations")
$wsMorning.Name = "pre-feeding-concentrations"

$wsAlgae = $wb.Worksheets.Item("algae-added")

# --- pre-feeding-concentrations (sheet1) : fill in row 4 ---
$wsMorning.Range("A4").Value = 42949
$wsMorning.Range("B4").Value = 8
$wsMorning.Range("C4").Value = 43
$wsMorning.Range("D4").Value = 28
$wsMorning.Range("E4").Value = 48
$wsMorning.Range("F4").Value = 33
$wsMorning.Range("G4").Value = 27
$wsMorning.Range("H4").Formula = "=AVERAGE(C4,D4,E4,F4,G4)"
$wsMorning.Range("I4").Formula = "=(H4*9)/0.0009"
$wsMorning.Range("J4").Formula = "=15000*I4"
$wsMorning.Range("K4").Formula = "=15000*50000"
$wsMorning.Range("L4").Formula = "=K4-J4"
$wsMorning.Range("H4:L4").NumberFormat = "0"
$wsMorning.Range("M4").Value = "Took this first thing in the morning, but I fed before I left yesterday. I'll feed in the afternoon today as well and start only feeding before I leave. I will need to start measuring food presence right before feeding in the afternoon."

$wsMorning.Range("M5").Select() | Out-Null

# --- algae-added (sheet2) : add row 5 ---
# Copy date-format (A) and wrap-text (B) styles down from row 4 so the new
# cells pick up the same cellXfs entries instead of synthesizing new ones.
$wsAlgae.Range("A4").Copy() | Out-Null
$wsAlgae.Range("A5").PasteSpecial(-4122) | Out-Null
$wsAlgae.Range("B4").Copy() | Out-Null
$wsAlgae.Range("B5").PasteSpecial(-4122) | Out-Null

$wsAlgae.Range("A5").Value = 42949
$wsAlgae.Range("B5").Value = "500 mL Ciso, 500 mL Chagra"
$wsAlgae.Range("C5").Value = 110
$wsAlgae.Range("D5").Value = 116
$wsAlgae.Range("E5").Value = 110
$wsAlgae.Range("F5").Value = 95
$wsAlgae.Range("G5").Value = 130
$wsAlgae.Range("H5").Formula = "=AVERAGE(C5:G5)"
$wsAlgae.Range("I5").Formula = "=(H5*9)/0.0009"
$wsAlgae.Range("J5").Formula = "=15000*50000"
$wsAlgae.Range("K5").Formula = "=J5/I5"
$wsAlgae.Range("L5").Value = 700
$wsAlgae.Range("M5").Formula = "=L5*I5"
$wsAlgae.Range("N5").Formula = "=M5/15000"

# Row 2 and row 4 (the other multi-line "Strains" rows) are 32pt tall; match it.
$wsAlgae.Rows.Item(5).RowHeight = 32

$wsAlgae.Range("L5").Select() | Out-Null
